# Update version string across the workbook for the
# "Coal Mine Boundaries and Methane Sources - version 1.0.0" release.

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$wb = $excel.ActiveWorkbook

$aboutSheet = $wb.Worksheets.Item("About")
$boundariesSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# --- About sheet ---

# A2: standalone version string
$aboutSheet.Range("A2").Value = "Version: " + $newVersion

# A6: Recommended citation containing the version string
$newCitation = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for La Ramada Coal Mine, Colombia, M6708, version ''' + $newVersion + '''. (See the CC license for attribution requirements if sharing or adapting the data set.)'
$aboutSheet.Range("A6").Value = $newCitation

# --- Boundaries and methane sources sheet ---
# Column S ("build_version") rows 2-8 hold the version string.
for ($row = 2; $row -le 8; $row++) {
    $cell = $boundariesSheet.Cells.Item($row, 19)
    $current = $cell.Value()
    if ($current -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
